$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mismatched "Traders! Super Session Support" labels so they correctly
# correspond to the "6" / "11" scenario rows.
$ws.Range("C3").Value = "Traders! Super Session Support 6"
$ws.Range("C4").Value = "Traders! Super Session Support 11"

# Move the active selection from the header row to C5.
$ws.Range("C5").Select()
